# Slide 2 ("Data 3550 - Programming" / intro+setup slide), Content Placeholder 2:
#   1) "...we’ll go at 12:10. " -> "...we’ll go at 12:05. " (splits into two runs)
#   2) Merge the two runs of the last bullet ("Enter the email used from " +
#      "GitHub “GitHub emails” quiz on the LMS. ") into a single run.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item("Content Placeholder 2")
$tr = $shape.TextFrame.TextRange

# --- Change 1: paragraph 1 - update the meeting time from 12:10 to 12:05 ---
# Paragraph 1 text: "You’ll need to do some stuff, please start, we’ll go at 12:10. "
# Characters 1-53  = "You’ll need to do some stuff, please start, we’ll go "
# Characters 54-63 = "at 12:10. "
$para1 = $tr.Paragraphs(1, 1)
$timePart = $para1.Characters(54, 10)
$timePart.Text = "at 12:05. "

# --- Change 2: last bullet - merge the two runs into a single run ---
# Paragraph text: "Enter the email used from " + "GitHub “GitHub emails” quiz on the LMS. "
$para7 = $tr.Paragraphs(7, 1)
$firstRun = $para7.Characters(1, 26)
$firstRun.Delete()
$remainder = $para7.Characters(1, 40)
$remainder.Text = "Enter the email used from GitHub “GitHub emails” quiz on the LMS. "
